# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    4  = 1613
    6  = 1095
    8  = 11482
    12 = 358
    15 = 12375
    16 = 13052
    18 = 141
    21 = 223
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
